$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("N2").Value = 0
$ws.Range("Q2").Value = 1
$ws.Range("T2").Value = 2

$ws.Range("N4").Value = 0.7041217312852566
$ws.Range("O4").Value = -1
$ws.Range("P4").Value = 0.4082434625705131
$ws.Range("Q4").Value = 0.5497199333020351
$ws.Range("R4").Value = 1
$ws.Range("S4").Value = 0.09943986660407012
$ws.Range("T4").Value = 0.7741065466434253
$ws.Range("U4").Value = 1
$ws.Range("V4").Value = 0.5482130932868505

$ws.Range("B15").Value = 9.869827596845477
$ws.Range("B16").Value = 4.951099843665786
$ws.Range("B17").Value = -9.458561131516781
$ws.Range("B18").Value = 6.552760643691096
$ws.Range("B19").Value = 9.988824899049497
$ws.Range("B20").Value = 13.28321426209553
$ws.Range("B21").Value = 23.70176620263562
